$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("NYKAA", "EF"),
    @("CAMS", "M"),
    @("CLEAN", "M"),
    @("RITES", "EP"),
    @("ZOMOTO", "M"),
    @("NUVOCO", "F"),
    @("LXCHEM", "F"),
    @("EASEMYTRIP-BE", "EP"),
    @("IPL", "F"),
    @("ROSSARI", "P"),
    @("UTIAMC", "M"),
    @("CHEMCON", "M"),
    @("HAPPSTMNDS-BE", "P"),
    @("IEX", "P"),
    @("CRAFTSMAN", "P")
)

$row = 4
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

$ws.Rows.Item(18).Select() | Out-Null
